$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Values are entered in the same order the original author typed them, so
# that new shared-string entries land in the expected index order.
$ws.Range("A23").Value = "Availability"
$ws.Range("B23").Value = "Online + Sourcecode"
$ws.Range("C23").Value = "Online + Sourcecode (old version)"
$ws.Range("D23").Value = "R"
$ws.Range("E23").Value = "R"
$ws.Range("F23").Value = "R"
$ws.Range("G23").Value = "bash pipeline/docker"
$ws.Range("A24").Value = "Runs"
$ws.Range("H23").Value = "Online + R package"
$ws.Range("H25").Value = "https://github.com/GfellerLab/EPIC"
$ws.Range("A25").Value = "Link"

$ws.Range("D24").Value = "x"
$ws.Range("E24").Value = "x"
$ws.Range("F24").Value = "x"
$ws.Range("H24").Value = "x"

# Bold style for column A labels, matching existing header cells in column A
$ws.Range("A23").Font.Bold = $true
$ws.Range("A24").Font.Bold = $true
$ws.Range("A25").Font.Bold = $true

# Column G width adjustment -> stored column width of 20. The engine adds a
# fixed 5/6-character padding on top of whatever is assigned to ColumnWidth
# (consistent with Excel's own char-width/pixel-padding conversion), so we
# back that padding out here to land on a stored width of exactly 20.
$ws.Range("G1").EntireColumn.ColumnWidth = 20 - 5/6

# Update selection to A25 to match final cursor position
$ws.Range("A25").Select()
